$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.548.49"
$ws.Range("E2").Value = "  -1.05%  "
$ws.Range("D3").Value = "3.826.50"
$ws.Range("E3").Value = "  +2.04%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'600.19"
$ws.Range("E5").Value = "  -0.36%  "
$ws.Range("D6").Value = "'163.44"
$ws.Range("E6").Value = "  -2.92%  "
$ws.Range("D7").Value = "3.825.28"
$ws.Range("E7").Value = "  +2.07%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -2.43%  "
$ws.Range("E10").Value = "  -3.99%  "
$ws.Range("E11").Value = "  -1.07%  "
$ws.Range("E12").Value = "  -0.79%  "
$ws.Range("D13").Value = "'36.77"
$ws.Range("E13").Value = "  -3.95%  "
$ws.Range("E14").Value = "  -2.82%  "
$ws.Range("D15").Value = "4.464.02"
$ws.Range("E15").Value = "  +2.04%  "
$ws.Range("D16").Value = "3.799.79"
$ws.Range("E16").Value = "  +1.06%  "
$ws.Range("D17").Value = "68.682.32"
$ws.Range("E17").Value = "  -0.82%  "
$ws.Range("E18").Value = "  +1.67%  "
$ws.Range("E19").Value = "  -0.57%  "
$ws.Range("D20").Value = "'17.05"
$ws.Range("E20").Value = "  -2.28%  "
$ws.Range("E21").Value = "  -1.27%  "
$ws.Range("D22").Value = "'483.89"
$ws.Range("E22").Value = "  -1.97%  "
$ws.Range("D23").Value = "'0.716"
$ws.Range("E23").Value = "  -1.84%  "
$ws.Range("D24").Value = "'0.0000159"
$ws.Range("E24").Value = "  +6.78%  "
$ws.Range("D25").Value = "'83.93"
$ws.Range("E25").Value = "  -1.15%  "
$ws.Range("D26").Value = "'2.23"
$ws.Range("E26").Value = "  -2.76%  "
$ws.Range("D27").Value = "'12.06"
$ws.Range("E27").Value = "  -2.12%  "
$ws.Range("B28").Value = "Dai"
$ws.Range("C28").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D28").Value = "'0.997"
$ws.Range("E28").Value = "  -0.26%  "
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").Value = "'9.98"
$ws.Range("E29").Value = "  -0.99%  "
$ws.Range("D30").Value = "'2.94"
$ws.Range("E30").Value = "  -1.59%  "
$ws.Range("B31").Value = "WrappedeETH"
$ws.Range("C31").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D31").Value = "3.974.86"
$ws.Range("E31").Value = "  +2.08%  "
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").Value = "'7.81"
$ws.Range("E32").Value = "  -4.85%  "
$ws.Range("E33").Value = "  -4.28%  "
$ws.Range("D34").Value = "'31.74"
$ws.Range("E34").Value = "  +0.20%  "
$ws.Range("D35").Value = "3.770.75"
$ws.Range("E35").Value = "  +2.41%  "
$ws.Range("E36").Value = "  -1.60%  "
$ws.Range("D37").Value = "'1.03"
$ws.Range("E37").Value = "  +1.02%  "
$ws.Range("E38").Value = "  -1.19%  "
$ws.Range("D39").Value = "'5.86"
$ws.Range("E39").Value = "  -2.57%  "
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("E41").Value = "  -3.26%  "
$ws.Range("D42").Value = "'434.13"
$ws.Range("E42").Value = "  +2.31%  "
$ws.Range("E43").Value = "  -4.27%  "
$ws.Range("D44").Value = "'48.46"
$ws.Range("E44").Value = "  -0.93%  "
$ws.Range("E45").Value = "  -1.00%  "
$ws.Range("E47").Value = "  -1.02%  "
$ws.Range("D48").Value = "2.835.54"
$ws.Range("E48").Value = "  +1.36%  "
$ws.Range("D49").Value = "'142.38"
$ws.Range("E49").Value = "  +0.84%  "
$ws.Range("E50").Value = "  +0.21%  "
$ws.Range("D51").Value = "'25.81"
$ws.Range("E51").Value = "  +12.74%  "
